# Edit script: updates the weekly Work Report for WR 89812901 to reflect the
# latest data pull: refreshed report timestamp, recalculated totals, cleared
# Scope ID, corrected Thursday (08/14) line items (one stale line item
# removed, causing everything below it to re-flow up by one row), and
# populated pricing for the Friday (08/15) line items + its TOTAL row, which
# also shifts up by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header / summary panel updates
# ---------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"
$ws.Range("C8").Value = 3292.74
$ws.Range("C9").Value = 8
$ws.Range("G10").Value = ""

# ---------------------------------------------------------------------
# Thursday (08/14/2025) table - rows 16-22
# Row 16 keeps its content; only the priced amount changes.
# ---------------------------------------------------------------------
$ws.Range("H16").Value = 54.9

# Row 17 now carries what used to be row 18's line item (priced).
$ws.Range("B17").Value = "CNC-HTA-10"
$ws.Range("D17").Value = "Compression connector H-Tap Assembly 1/0"
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 153.12

# Row 18 now carries what used to be row 19's line item (priced).
$ws.Range("B18").Value = "XFR-50-72-120-1B"
$ws.Range("D18").Value = "XFR,50KVA,7.2/12.4kVY,120/240,1BG"
$ws.Range("F18").Value = 1
$ws.Range("H18").Value = 234

# Row 19 now carries what used to be row 20's line item (priced).
$ws.Range("A19").Value = 60
$ws.Range("C19").Value = "Rem"
$ws.Range("H19").Value = 159

# Row 20 now carries what used to be row 21's line item (priced).
$ws.Range("B20").Value = "CNC-HTA-10"
$ws.Range("D20").Value = "Compression connector H-Tap Assembly 1/0"
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 97.59999999999999

# Row 21 now carries what used to be row 22's line item (priced).
$ws.Range("A21").Value = "Point 51"
$ws.Range("B21").Value = "PLA-HDIG"
$ws.Range("C21").Value = "Inst"
$ws.Range("D21").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("F21").Value = 2
$ws.Range("H21").Value = 1297.06

# Row 22 becomes the Thursday TOTAL row (it used to be the old row 23's
# spot). Pull its formatting from the old TOTAL row (row 23) first, then
# unmerge the old A23:G23 merge and set the new A22:G22 merge + values.
$ws.Range("A23:I23").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A23:G23").UnMerge()
$ws.Range("B22:G22").ClearContents()
$ws.Range("A22").Value = "TOTAL"
$ws.Range("H22").Value = 1995.68
$ws.Range("A22:G22").Merge()

# The old TOTAL row (23) is now fully vacated - clear it completely so it
# disappears from the sheet, matching the new dimension (...I29).
$ws.Range("A23:I23").Clear()

# ---------------------------------------------------------------------
# Friday (08/15/2025) table moves from rows 26-30 up to rows 25-29.
# Pull formatting down one row at a time (top to bottom, so each source
# row is copied before it is itself overwritten), then fill in values,
# then clear out the vacated old row 30.
# ---------------------------------------------------------------------
$ws.Range("A26:I26").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("A27:I27").Copy()
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("A28:I28").Copy()
$ws.Range("A27").PasteSpecial(-4122)

$ws.Range("A29:I29").Copy()
$ws.Range("A28").PasteSpecial(-4122)

$ws.Range("A30:I30").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Fix up the merged header/title bars for the new row positions.
$ws.Range("A26:H26").UnMerge()
$ws.Range("A30:G30").UnMerge()
$ws.Range("A25:H25").Merge()
$ws.Range("A29:G29").Merge()

# Row 25: section title.
$ws.Range("A25").Value = "Friday (08/15/2025)"

# Row 26: column headers (unchanged text, already correct from the copy).
$ws.Range("A26").Value = "Point Number"
$ws.Range("B26").Value = "Billable Unit Code"
$ws.Range("C26").Value = "Work Type"
$ws.Range("D26").Value = "Unit Description"
$ws.Range("E26").Value = "Unit of Measure"
$ws.Range("F26").Value = "# Units"
$ws.Range("G26").Value = "N/A"
$ws.Range("H26").Value = "Pricing"

# Row 27: Point 07 line item, now priced.
$ws.Range("A27").Value = "Point 07"
$ws.Range("B27").Value = "PLA-HDIG"
$ws.Range("C27").Value = "Inst"
$ws.Range("D27").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("E27").Value = "EA"
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 648.53

# Row 28: Point 11 line item, now priced.
$ws.Range("A28").Value = "Point 11"
$ws.Range("B28").Value = "PLA-HDIG"
$ws.Range("C28").Value = "Inst"
$ws.Range("D28").Value = "PLA,Hand Dig or Additional  Excavation"
$ws.Range("E28").Value = "EA"
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 648.53

# Row 29: Friday TOTAL row, now priced.
$ws.Range("B29:G29").ClearContents()
$ws.Range("A29").Value = "TOTAL"
$ws.Range("H29").Value = 1297.06

# The old Friday TOTAL row (30) is now fully vacated - clear it completely.
$ws.Range("A30:I30").Clear()
